$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "CubeA"

# Correct the slightly-off-precision values in row 15 (C15, G15)
$ws.Range("C15").Value = 1.25131244311135
$ws.Range("G15").Value = 1.25131244311135

# Append a new row 16 with Gaussian Quadrature averaged intensities data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.4592399385408076
$ws.Range("D16").Value = 2.343758461795015
$ws.Range("E16").Value = 1.226217356742332
$ws.Range("F16").Value = 0.8846105184970876
$ws.Range("G16").Value = 0.4592399385408076
$ws.Range("H16").Value = 2.343758461795015
$ws.Range("I16").Value = 0.9295949883073978
$ws.Range("J16").Value = 1.062061652702853
$ws.Range("K16").Value = 0.6296575682396935
$ws.Range("L16").Value = 1.396241255662679
$ws.Range("M16").Value = 0.4592399385408076
$ws.Range("N16").Value = 1.784987909268674
$ws.Range("O16").Value = 1.22845656889381
$ws.Range("P16").Value = 1.116422717560983

# Match the bold/bordered/centered style of A15 on the new A16 cell
# (copy formats only so the underlying cellXfs table is reused, not duplicated)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
